$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated forecast error values for Q1 (row2) through Q8 (row9)
$data = @{
    2 = @(0.0188948199481075, 0.4187755806181844, 0.4370514162648471, 0.6610986433693894, 0.6674041429934753, 51)
    3 = @(0.4263588555972532, 0.7006187285611748, 1.430786905557154, 1.196155050801172, 1.128935308341156, 50)
    4 = @(0.6005044672587662, 1.041385273230807, 4.110571973985722, 2.027454555344145, 1.956550704154995, 49)
    5 = @(0.4965930427113764, 1.09609697646386, 4.741065602765889, 2.177398815735392, 2.142449014907315, 48)
    6 = @(0.4101449545471543, 0.9701711184912666, 4.391446317909476, 2.095577800490708, 2.077266647837042, 47)
    7 = @(0.4037568587085331, 1.05621660122371, 5.365873945673021, 2.31643561224417, 2.314277311063799, 35)
    8 = @(0.3937173311530819, 1.056330459770304, 5.480532792941329, 2.341053778310385, 2.342413036496409, 34)
    9 = @(0.2639658233040533, 1.532043772134663, 10.17832809721954, 3.1903492124248, 3.277261169815604, 17)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
